$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.962.45"
$ws.Range("E2").Value = "  +3.63%  "
$ws.Range("D3").Value = "3.330.01"
$ws.Range("E3").Value = "  +1.09%  "
$ws.Range("E4").Value = "  -0.42%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  +3.95%  "
$ws.Range("D9").Value = "3.326.69"
$ws.Range("E9").Value = "  +1.28%  "
$ws.Range("E10").Value = "  +4.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.583"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.70"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000277"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "648.90"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +15.27%  "
$ws.Range("D15").Value = "3.857.64"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.79%  "
$ws.Range("D17").Value = "68.059.75"
$ws.Range("E17").Value = "  +3.85%  "
$ws.Range("E18").Value = "  +1.83%  "
$ws.Range("D19").Value = "3.326.57"
$ws.Range("E19").Value = "  +0.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.902"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.76"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.17%  "
$ws.Range("E24").Value = "  +2.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "98.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.14"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +10.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.60"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.67"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "602.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.39%  "
$ws.Range("D33").Value = "3.935.55"
$ws.Range("E33").Value = "  +4.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.99"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.52%  "
$ws.Range("E36").Value = "  +3.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.90"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.85%  "
$ws.Range("E39").Value = "  +6.38%  "
$ws.Range("E40").Value = "  +3.14%  "
$ws.Range("E41").Value = "  +7.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "33.04"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.29%  "
$ws.Range("D43").Value = "0.0₃0690"
$ws.Range("E43").Value = "  +3.88%  "
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.37"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.58%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.340"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0417"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.83%  "
$ws.Range("E47").Value = "  +3.44%  "
$ws.Range("B48").Value = "ThetaToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.57"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.12%  "
$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.35"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +10.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "130.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.34%  "
